$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the "Periodo Mora" column (E16:E29) to reverse-chronological
# order (most recent period first) and carry the corresponding "Valor Mora"
# values (F16:F29) along with them. The previous account-statement periods
# are being dropped from the top of the list and newer ones pushed in, per
# the commit message ("Elimina EC anteriores y se agregan nuevos").
$periods = @("2402","2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E" + $row).Value = $periods[$i]
}

# Value Mora for period 2402 is 33600 (all other periods keep 56000); since
# 2402 moved from the last data row (29) to the first (16), move its value
# with it and restore 56000 on the row that now holds period 2301.
$ws.Range("F16").Value = 33600
$ws.Range("F29").Value = 56000

# --- Column width tweaks (columns B, C, E, F, G, H, I, J got a bit wider)
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 9.333333333333334
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.333333333333332
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
